# Add debug message: populate the "Port-info" (column E) notes for a few
# ports, and rename the generic "busy"/"overrun_error"/"frame_error" wire
# names on the uart_rx / uart_tx sub-module sheets so they don't collide
# with the top-level "uart" instantiation sheet (which gets its own
# rx_busy/rx_frame_error/rx_overrun_error/tx_busy port rows appended).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "uart" (top-level instantiation) — drop the stand-alone "busy"
# row and append rx_busy / rx_frame_error / rx_overrun_error / tx_busy,
# plus two inline debug notes in column E.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Rows(15).Insert()
$ws1.Rows(15).RowHeight = 16

$ws1.Range("A3").Value = "clk"
$ws1.Range("B3").Value = "input"
$ws1.Range("C3").Value = 1

$ws1.Range("A4").Value = "m_axis_tready"
$ws1.Range("B4").Value = "input"
$ws1.Range("C4").Value = 1

$ws1.Range("A5").Value = "m_axis_tvalid"
$ws1.Range("B5").Value = "output"
$ws1.Range("C5").Value = 1

$ws1.Range("A6").Value = "prescale"
$ws1.Range("B6").Value = "input"
$ws1.Range("C6").Value = 16
$ws1.Range("E6").Value = "sfdaf"

$ws1.Range("A7").Value = "rst"
$ws1.Range("B7").Value = "input"
$ws1.Range("C7").Value = 1

$ws1.Range("A8").Value = "rxd"
$ws1.Range("B8").Value = "input"
$ws1.Range("C8").Value = 1

$ws1.Range("A9").Value = "s_axis_tready"
$ws1.Range("B9").Value = "output"
$ws1.Range("C9").Value = 1
$ws1.Range("E9").Value = "sdfdf"

$ws1.Range("A10").Value = "s_axis_tvalid"
$ws1.Range("B10").Value = "input"
$ws1.Range("C10").Value = 1

$ws1.Range("A11").Value = "txd"
$ws1.Range("B11").Value = "output"
$ws1.Range("C11").Value = 1

$ws1.Range("A12").Value = "rx_busy"
$ws1.Range("B12").Value = "output"
$ws1.Range("C12").Value = 1

$ws1.Range("A13").Value = "rx_frame_error"
$ws1.Range("B13").Value = "output"
$ws1.Range("C13").Value = 1

$ws1.Range("A14").Value = "rx_overrun_error"
$ws1.Range("B14").Value = "output"
$ws1.Range("C14").Value = 1

$ws1.Range("A15").Value = "tx_busy"
$ws1.Range("B15").Value = "output"
$ws1.Range("C15").Value = 1

# ---------------------------------------------------------------------
# Sheet "uart_rx" — rename the wire-name column entries for busy /
# overrun_error / frame_error so they match the renamed ports above, and
# drop in a few debug notes.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("E4").Value = "tyur"
$ws2.Range("E7").Value = "gfh"
$ws2.Range("D9").Value = "rx_busy"
$ws2.Range("D10").Value = "rx_overrun_error"
$ws2.Range("D11").Value = "rx_frame_error"
$ws2.Range("E11").Value = "hdgfh"

# ---------------------------------------------------------------------
# Sheet "uart_tx" — rename the wire-name column entry for busy, and add
# its debug notes.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("E4").Value = "sdfgfsdg"
$ws3.Range("D9").Value = "tx_busy"
$ws3.Range("E9").Value = "sfg"
